$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2: "USD" -> "EUR"
$ws.Range("A2").Value = "EUR"

# Update E2: 1000264788 -> 1000008617
$ws.Range("E2").Value = 1000008617

# Update selection to E2
$ws.Range("E2").Select()
